$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulated-probability matrix cells on Sheet1 with refreshed
# Monte Carlo results (more simulated games were added upstream).
$ws.Range("B2").Value = 0.1829652996845426
$ws.Range("C2").Value = 0.580441640378549
$ws.Range("J2").Value = 0.006309148264984227
$ws.Range("P2").Value = 0.1261829652996845
$ws.Range("S2").Value = 0.1041009463722398
$ws.Range("B3").Value = 0.01020408163265306
$ws.Range("C3").Value = 0.02040816326530612
$ws.Range("J3").Value = 0.03061224489795918
$ws.Range("P3").Value = 0.6989795918367347
$ws.Range("S3").Value = 0.2397959183673469
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("S4").Value = 0.3035714285714285
$ws.Range("B6").Value = 0.08866995073891626
$ws.Range("D6").Value = 0.01477832512315271
$ws.Range("F6").Value = 0.06403940886699508
$ws.Range("J6").Value = 0.2413793103448276
$ws.Range("O6").Value = 0.02463054187192118
$ws.Range("Q6").Value = 0.1428571428571428
$ws.Range("R6").Value = 0.04926108374384237
$ws.Range("S6").Value = 0.374384236453202
$ws.Range("B7").Value = 0.1160714285714286
$ws.Range("D7").Value = 0.04464285714285714
$ws.Range("F7").Value = 0.04464285714285714
$ws.Range("J7").Value = 0.15625
$ws.Range("O7").Value = 0.01339285714285714
$ws.Range("Q7").Value = 0.15625
$ws.Range("R7").Value = 0.06696428571428571
$ws.Range("S7").Value = 0.4017857142857143
$ws.Range("B8").Value = 0.1116625310173697
$ws.Range("D8").Value = 0.0173697270471464
$ws.Range("F8").Value = 0.05707196029776675
$ws.Range("J8").Value = 0.119106699751861
$ws.Range("O8").Value = 0.0173697270471464
$ws.Range("Q8").Value = 0.1488833746898263
$ws.Range("R8").Value = 0.06451612903225806
$ws.Range("S8").Value = 0.4640198511166253
$ws.Range("B9").Value = 0.08450704225352113
$ws.Range("D9").Value = 0.0352112676056338
$ws.Range("F9").Value = 0.07042253521126761
$ws.Range("J9").Value = 0.09859154929577464
$ws.Range("O9").Value = 0.0352112676056338
$ws.Range("Q9").Value = 0.1549295774647887
$ws.Range("R9").Value = 0.06338028169014084
$ws.Range("S9").Value = 0.4577464788732394
$ws.Range("B10").Value = 0.1215106732348112
$ws.Range("D10").Value = 0.02791461412151067
$ws.Range("F10").Value = 0.07142857142857142
$ws.Range("J10").Value = 0.1001642036124795
$ws.Range("O10").Value = 0.0180623973727422
$ws.Range("Q10").Value = 0.2019704433497537
$ws.Range("R10").Value = 0.06732348111658457
$ws.Range("S10").Value = 0.3916256157635468
$ws.Range("G11").Value = 0.1565217391304348
$ws.Range("J11").Value = 0.09565217391304348
$ws.Range("K11").Value = 0.2115942028985507
$ws.Range("L11").Value = 0.5101449275362319
$ws.Range("S11").Value = 0.02608695652173913
$ws.Range("G12").Value = 0.7197802197802198
$ws.Range("J12").Value = 0.2307692307692308
$ws.Range("K12").Value = 0.005494505494505495
$ws.Range("L12").Value = 0.01098901098901099
$ws.Range("S12").Value = 0.03296703296703297
$ws.Range("F13").Value = 0.0136986301369863
$ws.Range("G13").Value = 0.6027397260273972
$ws.Range("J13").Value = 0.2876712328767123
$ws.Range("S13").Value = 0.0958904109589041
$ws.Range("F15").Value = 0.03296703296703297
$ws.Range("H15").Value = 0.1153846153846154
$ws.Range("I15").Value = 0.07142857142857142
$ws.Range("J15").Value = 0.3736263736263736
$ws.Range("K15").Value = 0.06593406593406594
$ws.Range("M15").Value = 0.005494505494505495
$ws.Range("O15").Value = 0.03296703296703297
$ws.Range("S15").Value = 0.3021978021978022
$ws.Range("F16").Value = 0.03448275862068965
$ws.Range("H16").Value = 0.1724137931034483
$ws.Range("I16").Value = 0.06403940886699508
$ws.Range("J16").Value = 0.4137931034482759
$ws.Range("K16").Value = 0.1083743842364532
$ws.Range("M16").Value = 0.03940886699507389
$ws.Range("O16").Value = 0.02463054187192118
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.0178117048346056
$ws.Range("H17").Value = 0.1552162849872774
$ws.Range("I17").Value = 0.0916030534351145
$ws.Range("J17").Value = 0.4325699745547074
$ws.Range("K17").Value = 0.1246819338422392
$ws.Range("M17").Value = 0.0356234096692112
$ws.Range("O17").Value = 0.03816793893129771
$ws.Range("S17").Value = 0.1043256997455471
$ws.Range("F18").Value = 0.02127659574468085
$ws.Range("H18").Value = 0.1773049645390071
$ws.Range("I18").Value = 0.07801418439716312
$ws.Range("J18").Value = 0.4113475177304964
$ws.Range("K18").Value = 0.1276595744680851
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("O18").Value = 0.06382978723404255
$ws.Range("S18").Value = 0.09929078014184398
$ws.Range("F19").Value = 0.01098901098901099
$ws.Range("H19").Value = 0.2032967032967033
$ws.Range("I19").Value = 0.05259026687598116
$ws.Range("J19").Value = 0.3830455259026688
$ws.Range("K19").Value = 0.1326530612244898
$ws.Range("M19").Value = 0.03846153846153846
$ws.Range("O19").Value = 0.06200941915227629
$ws.Range("S19").Value = 0.1169544740973312

Write-Host "Updated 106 cells on $($ws.Name)"
